# Automatische test-sync: 2025-06-17 21:56:51
# Adds the newest incoming mail-log entry (row 37) to the "Logs" sheet,
# extends the conditional formatting on columns D and G to include the
# new row, and bumps the "Overig" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 37 -----------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A37").Value = "Vragen over samenwerking"
$logs.Range("B37").Value = "mailmind.test@zohomail.eu"
$logs.Range("C37").Value = "Kunnen we samenwerken aan een nieuw project?"
$logs.Range("D37").Value = "Overig"
$logs.Range("F37").Value = "2025-06-17 21:56:35"
$logs.Range("G37").Value = "Nee"

# Extend the existing conditional-formatting ranges so they keep
# covering the Categorie (D) and Beantwoord (G) columns through row 37.
$catFormat = $logs.Range("D2:D36").FormatConditions.Item(1)
$catFormat.ModifyAppliesToRange($logs.Range("D2:D37"))

$answeredFormat = $logs.Range("G2:G36").FormatConditions.Item(1)
$answeredFormat.ModifyAppliesToRange($logs.Range("G2:G37"))

# --- Dashboard sheet: bump the "Overig" count from 8 to 9 ------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B3").Value = 9
